$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "arrumar logo e finalizar paleta de cores"
# Finish filling out the color-palette table on the single worksheet.

# Row 20 (B20): fix the duplicated label "h2, h2" -> "h2, h3"
$ws.Range("B20").Value = "h2, h3"

# Row 31 (B31): extend the usage note with the final text reference
$ws.Range("B31").Value = "h2 dos cards, texto final"

# Row 35 (B35): replace the placeholder "Es" with the real usage label
$ws.Range("B35").Value = "h2 do blog"

# Row 35 (C35): this color cell was left blank - fill in the matching
# green hex code already used for "Titulo do blog" (row 7 / #78ef46)
$ws.Range("C35").Value = "#78ef46"

# Update the window's scroll position / active selection to match where
# the author ended up after finishing the edits.
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("C35").Select()
